# Applies the cryptocurrency price / 1h-volume refresh captured in the commit
# "Updated cryptos list on Sat Oct 19 23:51:32 UTC 2024 with GitHub Actions".
#
# All Coin/Link/Price/Volume(1h) cells in this sheet are stored as literal text
# (e.g. "68.352.55", "  -0.05%  ") rather than numbers, including values that
# look numeric. A plain `.Value = "598.25"` assignment would get auto-typed as a
# Number by Excel, so each cell is briefly switched to the Text number format,
# written, and then restored to its original style/format afterwards - exactly
# what happens when a user types a value into a cell that's formatted as Text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '68.361.97' },
    @{ Cell = "E2"; Value = '  -0.08%  ' },
    @{ Cell = "D3"; Value = '2.649.48' },
    @{ Cell = "E3"; Value = '  +0.22%  ' },
    @{ Cell = "D5"; Value = '598.25' },
    @{ Cell = "E5"; Value = '  -0.25%  ' },
    @{ Cell = "D6"; Value = '159.43' },
    @{ Cell = "E6"; Value = '  +2.81%  ' },
    @{ Cell = "E8"; Value = '  -0.39%  ' },
    @{ Cell = "E9"; Value = '  +5.25%  ' },
    @{ Cell = "E10"; Value = '  -0.78%  ' },
    @{ Cell = "D11"; Value = '5.27' },
    @{ Cell = "E11"; Value = '  +0.36%  ' },
    @{ Cell = "D12"; Value = '0.351' },
    @{ Cell = "E12"; Value = '  +0.23%  ' },
    @{ Cell = "D13"; Value = '28.08' },
    @{ Cell = "E13"; Value = '  -0.05%  ' },
    @{ Cell = "D14"; Value = '0.0000190' },
    @{ Cell = "E14"; Value = '  +1.24%  ' },
    @{ Cell = "D15"; Value = '3.132.31' },
    @{ Cell = "E15"; Value = '  +0.25%  ' },
    @{ Cell = "D16"; Value = '68.237.00' },
    @{ Cell = "E16"; Value = '  -0.03%  ' },
    @{ Cell = "D17"; Value = '2.665.11' },
    @{ Cell = "E17"; Value = '  +0.90%  ' },
    @{ Cell = "E18"; Value = '  -0.39%  ' },
    @{ Cell = "D19"; Value = '363.66' },
    @{ Cell = "E19"; Value = '  -0.93%  ' },
    @{ Cell = "E20"; Value = '  +3.04%  ' },
    @{ Cell = "E21"; Value = '  -1.77%  ' },
    @{ Cell = "E22"; Value = '  -1.06%  ' },
    @{ Cell = "E23"; Value = '  -3.15%  ' },
    @{ Cell = "D24"; Value = '75.09' },
    @{ Cell = "E25"; Value = '  +0.01%  ' },
    @{ Cell = "E26"; Value = '  -3.23%  ' },
    @{ Cell = "B27"; Value = 'Binance-PegBSC-USD' },
    @{ Cell = "C27"; Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd' },
    @{ Cell = "D27"; Value = '1.04' },
    @{ Cell = "E27"; Value = '  +3.73%  ' },
    @{ Cell = "B28"; Value = 'WrappedeETH' },
    @{ Cell = "C28"; Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth' },
    @{ Cell = "D28"; Value = '2.784.18' },
    @{ Cell = "E28"; Value = '  +0.41%  ' },
    @{ Cell = "E29"; Value = '  +0.04%  ' },
    @{ Cell = "D30"; Value = '558.35' },
    @{ Cell = "E30"; Value = '  -3.20%  ' },
    @{ Cell = "D31"; Value = '8.04' },
    @{ Cell = "E31"; Value = '  +0.22%  ' },
    @{ Cell = "E32"; Value = '  -1.65%  ' },
    @{ Cell = "E33"; Value = '  +0.32%  ' },
    @{ Cell = "E34"; Value = '  -0.87%  ' },
    @{ Cell = "E35"; Value = '  +0.03%  ' },
    @{ Cell = "D36"; Value = '1.58' },
    @{ Cell = "E36"; Value = '  +1.36%  ' },
    @{ Cell = "D37"; Value = '19.89' },
    @{ Cell = "E37"; Value = '  +2.80%  ' },
    @{ Cell = "D38"; Value = '159.72' },
    @{ Cell = "E38"; Value = '  -0.67%  ' },
    @{ Cell = "E39"; Value = '  +0.50%  ' },
    @{ Cell = "E40"; Value = '  -2.35%  ' },
    @{ Cell = "D41"; Value = '5.37' },
    @{ Cell = "E41"; Value = '  -0.85%  ' },
    @{ Cell = "E42"; Value = '  +4.66%  ' },
    @{ Cell = "B43"; Value = 'WhiteBITCoin' },
    @{ Cell = "C43"; Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt' },
    @{ Cell = "D43"; Value = '17.79' },
    @{ Cell = "E43"; Value = '  +0.28%  ' },
    @{ Cell = "B44"; Value = 'dogwifhat' },
    @{ Cell = "C44"; Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif' },
    @{ Cell = "D44"; Value = '2.64' },
    @{ Cell = "E44"; Value = '  -0.54%  ' },
    @{ Cell = "E45"; Value = '  +0.02%  ' },
    @{ Cell = "E46"; Value = '  -0.51%  ' },
    @{ Cell = "E47"; Value = '  +0.03%  ' },
    @{ Cell = "D48"; Value = '22.29' },
    @{ Cell = "E48"; Value = '  +1.00%  ' },
    @{ Cell = "E49"; Value = '  -1.56%  ' },
    @{ Cell = "D50"; Value = '0.0778' },
    @{ Cell = "E50"; Value = '  -0.26%  ' },
    @{ Cell = "D51"; Value = '0.616' },
    @{ Cell = "E51"; Value = '  -0.10%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $savedStyle = $cell.Style
    $cell.NumberFormat = "@"   # force text interpretation, like the cell already was
    $cell.Value = $u.Value
    $cell.Style = $savedStyle  # restore original (unformatted) style
}
